try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # A new survey record (row 34) is appended, duplicating row 32's
    # answers except for column H, which gets its own form-instance label.
    $sourceRow = 32
    $targetRow = 34

    $plainCols = @("A","B","C","D","E","F","G","J","K","L","M","N","O","P","Q")
    foreach ($col in $plainCols) {
        $ws.Range("$col$targetRow").Value = $ws.Range("$col$sourceRow").Value()
    }

    $ws.Range("H$targetRow").Value = "Evaluacin_Socioeconmica1580823537773 – 99"

    # Column I carries a zero-padded text code ("03"); force text so Excel
    # doesn't reinterpret it as the number 3.
    $ws.Range("I$targetRow").NumberFormat = "@"
    $ws.Range("I$targetRow").Value = $ws.Range("I$sourceRow").Value()
    $ws.Range("I$targetRow").NumberFormat = "General"

    # Widen column H to fit the new, longer label, and leave the active
    # selection on the next blank row like Excel does after data entry.
    $ws.Columns.Item(8).ColumnWidth = 40.14
    $ws.Range("H35").Select()
}
catch {
    throw
}
